$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.803.62"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.926.35"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'241.51"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4783"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("D8").Value = "'0.2884"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").Value = "'0.06789"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").Value = "'104.19"
$ws.Range("D12").Value = "'0.07793"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "1.933.10"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "'5.280"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "'0.6824"
$ws.Range("E15").Value = "  -2.89%  "
$ws.Range("D16").Value = "'291.84"
$ws.Range("E16").Value = "  +6.88%  "
$ws.Range("D17").Value = "30.799.64"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007583"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.185.74"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'12.88"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D22").Value = "'5.511"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("B23").Value = "BitDAO"
$ws.Range("C23").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D23").Value = "'0.4714"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'0.9998"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'6.385"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.538"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'168.07"
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.77"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.114"
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.391"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1008"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.604"
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.528"
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.326"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.04815"
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7352"
$ws.Range("E36").Value = "  -3.27%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.125"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'2.717"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01943"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.632"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.419"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'75.03"
$ws.Range("E42").Value = "  -5.74%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.022"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8681"
$ws.Range("E44").Value = "  -4.12%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4341"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'105.88"
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'0.9998"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.539"
$ws.Range("E48").Value = "  -4.20%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'987.87"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1212"
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").Value = "'9.021"
$ws.Range("E51").Value = "  -2.53%  "
